$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold text values that look numeric (trailing
# zeros, percents, very small magnitudes). Force Text format on each cell
# being rewritten so Excel keeps the literal string instead of silently
# reinterpreting it as a number (which would drop formatting like
# "3.500" -> 3.5 or "0.0001200" -> 1.2E-04).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "254.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.55%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.363"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.07%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.24%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.82%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8685"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9199"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7.51%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1421"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.89%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07166"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.99%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.17%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09249"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.37%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001566"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.64%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006054"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005813"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.33%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.500"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.22%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.232"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.222"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.01%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3180"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.66%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03461"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.60%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.48%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.538"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.36%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04171"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.23%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.05%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.005031"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "21.47%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001228"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.03%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001939"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "33.80%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03831"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1102"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.03%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003811"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-33.03%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002360"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01099"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "24.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005224"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.79%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08760"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "23.36%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002159"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.42%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
